$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from G1 ("sum") to H1,
# then set the new header text.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Save column values for rows 2-16 (1 = saved, 0 = not saved)
$saveValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 1
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 1
}

foreach ($row in $saveValues.Keys) {
    $ws.Cells.Item($row, 8).Value = $saveValues[$row]
}
